# Edit: add "2022-Q4" sheet with fund holding data, and add the corresponding
# summary row on the "总计" sheet (commit: feat: add 2022-Q4 data)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: shift existing quarters down one row
#    and insert the new 2022-Q4 totals at the top of the data (row 2).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summaryRows = @(
  @("2022-Q4", 39, 8.51),
  @("2022-Q3", 29, 8.36),
  @("2022-Q2", 18, 5.96),
  @("2022-Q1", 22, 6.91),
  @("2021-Q4", 24, 7.67),
  @("2021-Q3", 10, 2.9),
  @("2021-Q2", 2, 0.1),
  @("2021-Q1", 6, 2)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $row = $i + 2
    $entry = $summaryRows[$i]
    $summary.Cells.Item($row, 1).Value = $i
    $summary.Cells.Item($row, 2).Value = $entry[0]
    $summary.Cells.Item($row, 3).Value = $entry[1]
    $summary.Cells.Item($row, 4).Value = $entry[2]
}

# The newly appended row (old row 8, "2021-Q1") needs column A to carry the
# same bold/centered/bordered style ("s=2") used by every other index cell.
$summary.Range("A2").Copy()
$summary.Range("A9").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet right before "2022-Q3" (i.e. right
#    after "总计") and fill in the per-fund holding breakdown. Duplicating
#    the "2022-Q3" sheet (rather than Worksheets.Add) means the new tab
#    inherits the correct sheetPr / header+index styling / page setup, which
#    we then overwrite cell-by-cell with the 2022-Q4 figures.
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($q3Sheet)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# The template sheet only had 30 data rows; 2022-Q4 needs 39, so extend the
# bold index-column styling (column A) down through row 40.
$q4Sheet.Range("A2").Copy()
$q4Sheet.Range("A31:A40").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4Sheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$fundData = @"
0|012930|中庚价值先锋股票|68.71|94.78|4.46|3.0645|5
1|920003|中金新锐股票A|19.30|91.44|5.26|1.0152|5
2|007484|信澳核心科技混合|21.50|93.51|4.52|0.9718|3
3|410001|华富竞争力优选混合|10.06|83.31|4.90|0.4929|3
4|501078|广发科创主题灵活配置混合（LOF）|12.37|92.34|3.96|0.4899|6
5|009447|财通资管科技创新一年定期开放混合|9.13|94.71|4.60|0.4200|6
6|871003|广发资管价值增长灵活配置混合|6.51|88.29|4.37|0.2845|5
7|360007|光大保德信优势配置混合|7.90|75.56|3.09|0.2441|6
8|002064|华富产业升级灵活配置混合|6.57|79.23|3.59|0.2359|5
9|920923|中金新锐股票C|3.33|91.44|5.26|0.1752|5
10|162720|广发创业板两年定期开放混合|3.50|94.17|4.43|0.1550|7
11|020015|国泰区位优势混合A|2.05|84.45|5.92|0.1214|2
12|013067|富安达中小盘六个月持有期混合|2.09|78.87|4.97|0.1039|4
13|160642|鹏华增瑞灵活配置混合（LOF）|2.22|85.42|4.47|0.0992|6
14|015559|长江启航混合A|1.75|83.62|5.66|0.0990|5
15|163503|天治核心成长混合（LOF）|3.14|93.97|2.98|0.0936|10
16|580001|东吴嘉禾优势精选混合A|2.06|87.72|4.03|0.0830|8
17|005729|南方人工智能主题混合|2.01|89.09|3.83|0.0770|9
18|002577|南方新兴龙头灵活配置混合|1.42|86.40|3.92|0.0557|7
19|015026|鹏华增华混合A|1.65|78.57|2.94|0.0485|10
20|015594|国泰区位优势混合C|0.75|84.45|5.92|0.0444|2
21|012669|南方新兴产业混合A|0.66|88.33|3.92|0.0259|5
22|012670|南方新兴产业混合C|0.38|88.33|3.92|0.0149|5
23|011214|招商惠润一年定期开放混合（MOM）A|0.48|68.20|3.04|0.0146|6
24|350002|天治低碳经济灵活配置混合|0.65|85.45|1.81|0.0118|1
25|004608|长信乐信灵活配置混合A|0.49|32.30|2.06|0.0101|3
26|001709|华富物联世界灵活配置混合|0.20|88.45|4.92|0.0098|2
27|009128|明亚价值长青混合A|0.39|52.21|2.49|0.0097|7
28|350009|天治研究驱动混合A|0.29|92.96|3.18|0.0092|4
29|015027|鹏华增华混合C|0.20|78.57|2.94|0.0059|10
30|002043|天治研究驱动混合C|0.17|92.96|3.18|0.0054|4
31|002303|金鹰智慧生活灵活配置混合|0.09|94.65|5.17|0.0047|5
32|015560|长江启航混合C|0.07|83.62|5.66|0.0040|5
33|002584|富安达长盈灵活配置混合A|0.10|85.18|3.58|0.0036|5
34|011215|招商惠润一年定期开放混合（MOM）C|0.06|68.20|3.04|0.0018|6
35|004609|长信乐信灵活配置混合C|0.03|32.30|2.06|0.0006|3
36|015152|东吴嘉禾优势精选混合C|0.01|87.72|4.03|0.0004|8
37|016214|富安达长盈灵活配置混合C|0.01|85.18|3.58|0.0004|5
38|009129|明亚价值长青混合C|0.00|52.21|2.49|0|7
"@

$lines = $fundData -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $row = [int]$parts[0] + 2

    $q4Sheet.Cells.Item($row, 1).Value = [int]$parts[0]

    $q4Sheet.Cells.Item($row, 2).Value = "'" + $parts[1]
    $q4Sheet.Cells.Item($row, 2).Style = "Normal"

    $q4Sheet.Cells.Item($row, 3).Value = $parts[2]

    $q4Sheet.Cells.Item($row, 4).Value = "'" + $parts[3]
    $q4Sheet.Cells.Item($row, 4).Style = "Normal"

    $q4Sheet.Cells.Item($row, 5).Value = "'" + $parts[4]
    $q4Sheet.Cells.Item($row, 5).Style = "Normal"

    $q4Sheet.Cells.Item($row, 6).Value = "'" + $parts[5]
    $q4Sheet.Cells.Item($row, 6).Style = "Normal"

    $mktVal = $parts[6]
    if ($mktVal -eq "0") {
        $q4Sheet.Cells.Item($row, 7).Value = 0
    } else {
        $q4Sheet.Cells.Item($row, 7).Value = "'" + $mktVal
        $q4Sheet.Cells.Item($row, 7).Style = "Normal"
    }

    $q4Sheet.Cells.Item($row, 8).Value = [int]$parts[7]
}
